$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values to match the re-pulled data / mean calculation
$ws.Range("F4").Value = 22
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = -8
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -6
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 2
